$d = $word.ActiveDocument

# The first (and only) inline picture in the document sits alone in the
# "FirstParagraph"-styled paragraph right after the "Building Setback from
# Boundary" heading. Replace it with a hyperlink run pointing at the
# image's original URL on ura.gov.sg, styled with the built-in Hyperlink
# character style - matching how the rest of the document already links
# to other bookmarks/URLs.

$url = "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Industrial/BP02_Road_Buffer_Setbacks.jpg?h=100%25&w=100%25"

$shp = $d.InlineShapes.Item(1)
$rng = $shp.Range
$shp.Delete()

$d.Hyperlinks.Add($rng, $url, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $url) | Out-Null

Write-Output "Replaced inline picture with hyperlink to $url"
